$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the Price/Volume columns so purely numeric-looking
# strings (e.g. "0.9996") are NOT auto-converted to floating point values -
# matches the source inlineStr cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.628.20"
$ws.Range("D3").Value = "1.872.97"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "247.85"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.4733"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "0.2917"
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("D9").Value = "0.06481"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "22.07"
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("D11").Value = "0.07699"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").Value = "96.61"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").Value = "0.7379"
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("D14").Value = "1.867.11"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "5.151"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "273.14"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "30.604.47"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "13.31"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "0.000007517"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "2.114.93"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "0.9989"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "5.256"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "6.174"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "9.217"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").Value = "163.91"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").Value = "18.78"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "1.911"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "0.09980"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").Value = "1.345"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").Value = "1.510"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "4.282"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "4.100"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").Value = "0.04782"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "1.120"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").Value = "0.6956"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "0.01847"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "2.753"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").Value = "6.200"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").Value = "73.17"
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("D43").Value = "0.4178"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "0.8326"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("D46").Value = "101.72"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "9.352"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").Value = "35.44"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "6.979"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").Value = "917.96"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "0.05648"
$ws.Range("E51").Value = "  +1.29%  "

# Restore default cell style (style index 0) so no stray "s" attribute is
# written for these cells - matches the original workbook formatting.
$ws.Range("D2:E51").Style = "Normal"
